$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2013.8182
$ws.Range("I9").Value = 437.16666
$ws.Range("K9").Value = 437.16666
$ws.Range("M9").Value = -268.16666
$ws.Range("H28").Value = 1125.625
$ws.Range("I28").Value = 1066.5
$ws.Range("J28").Value = 1303
$ws.Range("K28").Value = 1066.5
$ws.Range("L28").Value = 1303
$ws.Range("M28").Value = -581.5
$ws.Range("N28").Value = -2273
$ws.Range("H64").Value = 11988.895
$ws.Range("I64").Value = 7649.5835
$ws.Range("K64").Value = 7649.5835
$ws.Range("M64").Value = -7401.5835
$ws.Range("H67").Value = 11988.895
$ws.Range("I67").Value = 7649.5835
$ws.Range("K67").Value = 7649.5835
$ws.Range("M67").Value = -6791.5835
$ws.Range("H107").Value = 662.61536
$ws.Range("I107").Value = 526.3
$ws.Range("K107").Value = 526.3
$ws.Range("M107").Value = 1393.7
$ws.Range("H113").Value = 5462.25
$ws.Range("I113").Value = 3124.5
$ws.Range("J113").Value = 7800
$ws.Range("K113").Value = 3124.5
$ws.Range("L113").Value = 7800
$ws.Range("M113").Value = 129.5
$ws.Range("N113").Value = -14308
$ws.Range("H132").Value = 3340.7646
$ws.Range("I132").Value = 2029.3
$ws.Range("J132").Value = 5214.2856
$ws.Range("K132").Value = 6087.9
$ws.Range("L132").Value = 15642.8568
$ws.Range("M132").Value = -3557.9
$ws.Range("N132").Value = -20702.8568
$ws.Range("H135").Value = 2063
$ws.Range("I135").Value = 2069.75
$ws.Range("J135").Value = 2036
$ws.Range("K135").Value = 18627.75
$ws.Range("L135").Value = 18324
$ws.Range("M135").Value = -16092.75
$ws.Range("N135").Value = -23394

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 457
$ws.Range("J5").Value = 199
$ws.Range("L5").Value = 199
$ws.Range("N5").Value = -423
$ws.Range("H32").Value = 4445.7803
$ws.Range("I32").Value = 3806.925
$ws.Range("K32").Value = 3806.925
$ws.Range("M32").Value = -3519.925
$ws.Range("H102").Value = 2040.4166
$ws.Range("I102").Value = 1248.5
$ws.Range("K102").Value = 1248.5
$ws.Range("M102").Value = 373.5
$ws.Range("H110").Value = 2511.6428
$ws.Range("I110").Value = 1464.8889
$ws.Range("J110").Value = 4395.8
$ws.Range("K110").Value = 1464.8889
$ws.Range("L110").Value = 4395.8
$ws.Range("M110").Value = 580.1111000000001
$ws.Range("N110").Value = -8485.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 457
$ws.Range("J4").Value = 199
$ws.Range("L4").Value = 199
$ws.Range("N4").Value = -429
$ws.Range("H16").Value = 690
$ws.Range("I16").Value = 562.5
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 562.5
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -392.5
$ws.Range("N16").Value = -1540
$ws.Range("H22").Value = 796.125
$ws.Range("I22").Value = 811.6667
$ws.Range("J22").Value = 749.5
$ws.Range("K22").Value = 811.6667
$ws.Range("L22").Value = 749.5
$ws.Range("M22").Value = -638.6667
$ws.Range("N22").Value = -1095.5
$ws.Range("H36").Value = 1409
$ws.Range("I36").Value = 212
$ws.Range("J36").Value = 5000
$ws.Range("K36").Value = 212
$ws.Range("L36").Value = 5000
$ws.Range("M36").Value = 322
$ws.Range("N36").Value = -6068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1827
$ws.Range("I22").Value = 1058
$ws.Range("K22").Value = 1058
$ws.Range("M22").Value = -708
$ws.Range("H134").Value = 2429.9678
$ws.Range("I134").Value = 2230.8147
$ws.Range("K134").Value = 6692.4441
$ws.Range("M134").Value = -4157.4441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 111377890
$ws.Range("I4").Value = 300126
$ws.Range("K4").Value = 900378
$ws.Range("M4").Value = -900266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4397.8667
$ws.Range("I122").Value = 4499.75
$ws.Range("K122").Value = 13499.25
$ws.Range("M122").Value = -11049.25
$ws.Range("H126").Value = 200001150
$ws.Range("I126").Value = 200001150
$ws.Range("K126").Value = 600003450
$ws.Range("M126").Value = -600000980

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 251307.38
$ws.Range("I10").Value = 1000322.5
$ws.Range("J10").Value = 1635.6666
$ws.Range("K10").Value = 1000322.5
$ws.Range("L10").Value = 1635.6666
$ws.Range("M10").Value = -1000182.5
$ws.Range("N10").Value = -1915.6666
$ws.Range("H19").Value = 4317.875
$ws.Range("J19").Value = 8075.5
$ws.Range("L19").Value = 8075.5
$ws.Range("N19").Value = -8415.5
$ws.Range("H22").Value = 4285.7144
$ws.Range("I22").Value = 4000
$ws.Range("J22").Value = 4333.3335
$ws.Range("K22").Value = 4000
$ws.Range("L22").Value = 4333.3335
$ws.Range("M22").Value = -3705
$ws.Range("N22").Value = -4923.3335
$ws.Range("H27").Value = 4285.7144
$ws.Range("I27").Value = 4000
$ws.Range("J27").Value = 4333.3335
$ws.Range("K27").Value = 4000
$ws.Range("L27").Value = 4333.3335
$ws.Range("M27").Value = -3893
$ws.Range("N27").Value = -4547.3335
$ws.Range("H40").Value = 11756.632
$ws.Range("J40").Value = 9856.857
$ws.Range("L40").Value = 9856.857
$ws.Range("N40").Value = -10128.857
$ws.Range("H46").Value = 1532.4166
$ws.Range("I46").Value = 1900
$ws.Range("J46").Value = 1458.9
$ws.Range("K46").Value = 1900
$ws.Range("L46").Value = 1458.9
$ws.Range("M46").Value = -1712
$ws.Range("N46").Value = -1834.9
$ws.Range("H68").Value = 950.2
$ws.Range("I68").Value = 950.2
$ws.Range("K68").Value = 950.2
$ws.Range("M68").Value = -201.2
$ws.Range("H71").Value = 950.2
$ws.Range("I71").Value = 950.2
$ws.Range("K71").Value = 4751
$ws.Range("M71").Value = -1007
$ws.Range("H74").Value = 33999.8
$ws.Range("J74").Value = 36249.75
$ws.Range("L74").Value = 36249.75
$ws.Range("N74").Value = -38245.75
$ws.Range("H77").Value = 33999.8
$ws.Range("J77").Value = 36249.75
$ws.Range("L77").Value = 108749.25
$ws.Range("N77").Value = -118733.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H126").Value = 1719.7858
$ws.Range("I126").Value = 1422.25
$ws.Range("K126").Value = 4266.75
$ws.Range("M126").Value = -1796.75
$ws.Range("H132").Value = 2328.5
$ws.Range("I132").Value = 2192.5217
$ws.Range("J132").Value = 5456
$ws.Range("K132").Value = 6577.5651
$ws.Range("L132").Value = 16368
$ws.Range("M132").Value = -4047.5651
$ws.Range("N132").Value = -21428
